$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new price values parse as plain numbers need to be forced
# back to Text format first, otherwise Excel auto-converts the literal
# string (e.g. "301.70") into a numeric value (301.7) and the cell loses
# its original text/shared-string nature, which the source data relies on.
$textForcedCells = @(
    "D5", "D6", "D9", "D10", "D11", "D13", "D17", "D19",
    "D20", "D22", "D23", "D27", "D29", "D32", "D33", "D34",
    "D37", "D40", "D42", "D44", "D46", "D47", "D50", "D51"
)
foreach ($cellRef in $textForcedCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range("D2").Value = "42.897.63"
$ws.Range("E2").Value = "  +0.07%  "
$ws.Range("D3").Value = "2.368.87"
$ws.Range("E3").Value = "  +2.16%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "301.70"
$ws.Range("E5").Value = "  -0.23%  "
$ws.Range("D6").Value = "96.03"
$ws.Range("E6").Value = "  -0.21%  "
$ws.Range("E7").Value = "  -0.54%  "
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("D9").Value = "0.496"
$ws.Range("E9").Value = "  +0.05%  "
$ws.Range("D10").Value = "34.14"
$ws.Range("D11").Value = "0.0788"
$ws.Range("E11").Value = "  +0.42%  "
$ws.Range("E12").Value = "  +2.59%  "
$ws.Range("D13").Value = "18.21"
$ws.Range("E13").Value = "  -3.81%  "
$ws.Range("E14").Value = "  -0.18%  "
$ws.Range("D15").Value = "2.742.57"
$ws.Range("E15").Value = "  +2.17%  "
$ws.Range("D16").Value = "2.378.18"
$ws.Range("E16").Value = "  +2.82%  "
$ws.Range("D17").Value = "0.801"
$ws.Range("E17").Value = "  +1.30%  "
$ws.Range("D18").Value = "42.892.41"
$ws.Range("E18").Value = "  +0.20%  "
$ws.Range("D19").Value = "12.15"
$ws.Range("E19").Value = "  -0.33%  "
$ws.Range("D20").Value = "6.31"
$ws.Range("E20").Value = "  +2.41%  "
$ws.Range("D21").Value = "0.0₃0888"
$ws.Range("E21").Value = "  -0.60%  "
$ws.Range("D22").Value = "68.01"
$ws.Range("E22").Value = "  +0.08%  "
$ws.Range("D23").Value = "234.91"
$ws.Range("E24").Value = "  -1.58%  "
$ws.Range("E25").Value = "  -0.03%  "
$ws.Range("E26").Value = "  +0.73%  "
$ws.Range("D27").Value = "24.83"
$ws.Range("E27").Value = "  +1.79%  "
$ws.Range("E28").Value = "  +0.30%  "
$ws.Range("D29").Value = "9.23"
$ws.Range("E29").Value = "  +0.99%  "
$ws.Range("E30").Value = "  -2.27%  "
$ws.Range("E31").Value = "  +0.03%  "
$ws.Range("D32").Value = "5.05"
$ws.Range("E32").Value = "  +0.93%  "
$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D33").Value = "0.0733"
$ws.Range("E33").Value = "  +4.62%  "
$ws.Range("B34").Value = "Celestia"
$ws.Range("C34").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D34").Value = "17.42"
$ws.Range("E34").Value = "  -2.55%  "
$ws.Range("E35").Value = "  +5.28%  "
$ws.Range("E36").Value = "  +3.97%  "
$ws.Range("D37").Value = "4.35"
$ws.Range("E37").Value = "  -2.67%  "
$ws.Range("E38").Value = "  -0.71%  "
$ws.Range("E39").Value = "  +1.63%  "
$ws.Range("D40").Value = "22.22"
$ws.Range("E40").Value = "  +7.07%  "
$ws.Range("E41").Value = "  -0.65%  "
$ws.Range("D42").Value = "116.29"
$ws.Range("E42").Value = "  -30.11%  "
$ws.Range("D43").Value = "1.940.47"
$ws.Range("E43").Value = "  +0.25%  "
$ws.Range("D44").Value = "0.0280"
$ws.Range("E44").Value = "  +0.29%  "
$ws.Range("E45").Value = "  +1.91%  "
$ws.Range("D46").Value = "2.73"
$ws.Range("E46").Value = "  -1.16%  "
$ws.Range("D47").Value = "9.19"
$ws.Range("E47").Value = "  -9.87%  "
$ws.Range("D48").Value = "2.602.97"
$ws.Range("E48").Value = "  +2.01%  "
$ws.Range("E49").Value = "  +1.88%  "
$ws.Range("B50").Value = "BitcoinSV"
$ws.Range("C50").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D50").Value = "72.02"
$ws.Range("E50").Value = "  -0.13%  "
$ws.Range("B51").Value = "MultiversX"
$ws.Range("C51").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D51").Value = "52.02"
$ws.Range("E51").Value = "  -2.62%  "
